$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update scoring values in the module3 schedule block (columns E-H, rows 2-29)
$ws.Range("H2").Value  = "5-20"
$ws.Range("G3").Value  = "0-5"
$ws.Range("G4").Value  = "5-20"
$ws.Range("E5").Value  = "0-5"
$ws.Range("E6").Value  = "0-0"
$ws.Range("G7").Value  = "5-5"
$ws.Range("F8").Value  = "20-20"
$ws.Range("E9").Value  = "5-10"
$ws.Range("F10").Value = "20-20"
$ws.Range("E11").Value = "5-10"
$ws.Range("F11").Value = "5-2"
$ws.Range("G12").Value = "0-0"
$ws.Range("F13").Value = "5-5"
$ws.Range("H14").Value = "5-10"
$ws.Range("F15").Value = "10-5"
$ws.Range("G15").Value = "0-5"
$ws.Range("G16").Value = "5-20"
$ws.Range("E17").Value = "10-25"
$ws.Range("G19").Value = "5-10"
$ws.Range("E20").Value = "10-20"
$ws.Range("G21").Value = "0-0"
$ws.Range("E22").Value = "10-25"
$ws.Range("E23").Value = "5-20"
$ws.Range("H24").Value = "10-5"
$ws.Range("E25").Value = "5-10"
$ws.Range("G26").Value = "0-0"
$ws.Range("H27").Value = "10-20"
$ws.Range("E28").Value = "20-0"
$ws.Range("F28").Value = "5-10"
$ws.Range("E29").Value = "20-0"

# sheet view: update the active selection to reflect where the author left off
$ws.Range("H27").Select()

$wb.Save()
